$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 6502226.5
$ws.Range("I106").Value = 8822094
$ws.Range("K106").Value = 8822094
$ws.Range("M106").Value = -8821463
$ws.Range("H131").Value = 7778.0713
$ws.Range("J131").Value = 9296.591
$ws.Range("L131").Value = 27889.773
$ws.Range("N131").Value = -37969.773
$ws.Range("H138").Value = 1520.8904
$ws.Range("I138").Value = 847.05457
$ws.Range("J138").Value = 3579.8333
$ws.Range("K138").Value = 2541.16371
$ws.Range("L138").Value = 10739.4999
$ws.Range("M138").Value = 2598.83629
$ws.Range("N138").Value = -21019.4999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6794
$ws.Range("I32").Value = 6843.8535
$ws.Range("K32").Value = 6843.8535
$ws.Range("M32").Value = -6556.8535
$ws.Range("H45").Value = 8833.223
$ws.Range("I45").Value = 8874.75
$ws.Range("K45").Value = 8874.75
$ws.Range("M45").Value = -8497.75
$ws.Range("H61").Value = 3573.8916
$ws.Range("I61").Value = 3521.3796
$ws.Range("K61").Value = 3521.3796
$ws.Range("M61").Value = -3309.3796
$ws.Range("H74").Value = 7929.593
$ws.Range("I74").Value = 8283.556
$ws.Range("J74").Value = 7221.6665
$ws.Range("K74").Value = 8283.556
$ws.Range("L74").Value = 7221.6665
$ws.Range("M74").Value = -7409.556
$ws.Range("N74").Value = -8969.666499999999
$ws.Range("H77").Value = 7929.593
$ws.Range("I77").Value = 8283.556
$ws.Range("J77").Value = 7221.6665
$ws.Range("K77").Value = 41417.78
$ws.Range("L77").Value = 36108.3325
$ws.Range("M77").Value = -37049.78
$ws.Range("N77").Value = -44844.3325
$ws.Range("H122").Value = 2322530
$ws.Range("I122").Value = 7322
$ws.Range("K122").Value = 21966
$ws.Range("M122").Value = -19516
$ws.Range("H136").Value = 3573.8916
$ws.Range("I136").Value = 3521.3796
$ws.Range("K136").Value = 10564.1388
$ws.Range("M136").Value = -8014.138800000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 585
$ws.Range("I10").Value = 706.25
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 706.25
$ws.Range("L10").Value = 100
$ws.Range("M10").Value = -566.25
$ws.Range("N10").Value = -380
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3195.08
$ws.Range("I31").Value = 1996.4
$ws.Range("K31").Value = 1996.4
$ws.Range("M31").Value = -1701.4
$ws.Range("H34").Value = 3195.08
$ws.Range("I34").Value = 1996.4
$ws.Range("K34").Value = 1996.4
$ws.Range("M34").Value = -1794.4
$ws.Range("H58").Value = 2154.92
$ws.Range("I58").Value = 1038.5333
$ws.Range("K58").Value = 1038.5333
$ws.Range("M58").Value = -835.5333000000001
$ws.Range("H94").Value = 2762.077
$ws.Range("J94").Value = 3142.5557
$ws.Range("L94").Value = 3142.5557
$ws.Range("N94").Value = -4044.5557
$ws.Range("H119").Value = 94666.664
$ws.Range("I119").Value = 90000
$ws.Range("J119").Value = 97000
$ws.Range("K119").Value = 90000
$ws.Range("L119").Value = 97000
$ws.Range("M119").Value = -85162
$ws.Range("N119").Value = -106676
$ws.Range("H132").Value = 19641.75
$ws.Range("I132").Value = 2597.25
$ws.Range("K132").Value = 7791.75
$ws.Range("M132").Value = -5261.75
$ws.Range("H136").Value = 2154.92
$ws.Range("I136").Value = 1038.5333
$ws.Range("K136").Value = 3115.5999
$ws.Range("M136").Value = -565.5999000000002
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 455836.97
$ws.Range("J5").Value = 771076.6
$ws.Range("L5").Value = 2313229.8
$ws.Range("N5").Value = -2313453.8
$ws.Range("H17").Value = 1375.4762
$ws.Range("I17").Value = 324.75
$ws.Range("J17").Value = 1622.7059
$ws.Range("K17").Value = 974.25
$ws.Range("L17").Value = 4868.1177
$ws.Range("M17").Value = -805.25
$ws.Range("N17").Value = -5206.1177
$ws.Range("H23").Value = 5952743
$ws.Range("J23").Value = 8333670.5
$ws.Range("L23").Value = 25001011.5
$ws.Range("N23").Value = -25001481.5
$ws.Range("H113").Value = 4537.9688
$ws.Range("J113").Value = 4909.3105
$ws.Range("L113").Value = 14727.9315
$ws.Range("N113").Value = -19067.9315
$ws.Range("H132").Value = 101379
$ws.Range("I132").Value = 898
$ws.Range("J132").Value = 126499.25
$ws.Range("K132").Value = 8082
$ws.Range("L132").Value = 1138493.25
$ws.Range("M132").Value = -5552
$ws.Range("N132").Value = -1143553.25
$ws.Range("H135").Value = 455836.97
$ws.Range("J135").Value = 771076.6
$ws.Range("L135").Value = 6939689.399999999
$ws.Range("N135").Value = -6944759.399999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2231.6667
$ws.Range("I132").Value = 2253.25
$ws.Range("J132").Value = 1800
$ws.Range("K132").Value = 6759.75
$ws.Range("L132").Value = 5400
$ws.Range("M132").Value = -4229.75
$ws.Range("N132").Value = -10460
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 60000
$ws.Range("J36").Value = 60000
$ws.Range("L36").Value = 60000
$ws.Range("N36").Value = -61124
$ws.Range("H40").Value = 25637.182
$ws.Range("I40").Value = 43298.7
$ws.Range("J40").Value = 10919.25
$ws.Range("K40").Value = 43298.7
$ws.Range("L40").Value = 10919.25
$ws.Range("M40").Value = -43162.7
$ws.Range("N40").Value = -11191.25
$ws.Range("H93").Value = 10618.6
$ws.Range("J93").Value = 5997.75
$ws.Range("L93").Value = 5997.75
$ws.Range("N93").Value = -8493.75
$ws.Range("H100").Value = 2005.4615
$ws.Range("I100").Value = 1897.3636
$ws.Range("K100").Value = 1897.3636
$ws.Range("M100").Value = -1356.3636
$ws.Range("H122").Value = 8527.736999999999
$ws.Range("I122").Value = 9974.583000000001
$ws.Range("J122").Value = 6047.4287
$ws.Range("K122").Value = 29923.749
$ws.Range("L122").Value = 18142.2861
$ws.Range("M122").Value = -27473.749
$ws.Range("N122").Value = -23042.2861
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 20510
$ws.Range("I81").Value = 43397.75
$ws.Range("J81").Value = 2199.8
$ws.Range("K81").Value = 86795.5
$ws.Range("L81").Value = 4399.6
$ws.Range("M81").Value = -85734.5
$ws.Range("N81").Value = -6521.6
$ws.Range("H84").Value = 20510
$ws.Range("I84").Value = 43397.75
$ws.Range("J84").Value = 2199.8
$ws.Range("K84").Value = 433977.5
$ws.Range("L84").Value = 21998
$ws.Range("M84").Value = -428673.5
$ws.Range("N84").Value = -32606
$ws.Range("H122").Value = 25415.32
$ws.Range("I122").Value = 3198.6667
$ws.Range("K122").Value = 9596.000100000001
$ws.Range("M122").Value = -7146.000100000001
